# Regenerate save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals (column G, header "K").
#
# Only the "K" column (column G) values change for rows 2-39; every other
# column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 3
    6  = 2
    7  = 2
    8  = 1
    9  = 3
    10 = 2
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 2
    16 = 2
    17 = 0
    18 = 3
    19 = 3
    20 = 1
    21 = 0
    22 = 1
    23 = 3
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 1
    32 = 3
    33 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
